$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row for "brain" (row 5), shifting subsequent rows up.
$ws.Rows.Item(5).Delete()
